$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Multiply the percent_recovery values in C2:C11 by 100
# (adjust percent recovery calculation per commit message)
$range = $ws.Range("C2:C11")
foreach ($cell in $range.Cells) {
    $old = $cell.Value()
    $cell.Value = $old * 100
}
